# courrier_radiation.template.docx - split a few single runs into
# multiple runs (identical formatting) and drop a stray trailing space.
#
# Word's Find/Replace keeps a hit inside its existing run when the
# replacement text has the same formatting, so to get genuinely separate
# <w:r> elements at the split points we touch (and immediately restore)
# a formatting property on the tail sub-range; Word then materialises a
# fresh run at that boundary even though the final formatting matches.

$d = $word.ActiveDocument
$nbsp = [char]0xA0

function Split-RunAt($range, $offset) {
    # Force a run boundary at $range.Start + $offset by nudging a
    # character property on the tail and putting it back.
    $tailStart = $range.Start + $offset
    $tailEnd = $range.End
    $tail = $d.Range($tailStart, $tailEnd)
    $tail.Font.Size = 12
    $tail.Font.Size = 11
}

# --- Change 1 -----------------------------------------------------
# "Objet<nbsp>: Résiliation de l'élection de domicile" (one run) becomes
# three runs: "Objet" | "<nbsp>: " | "Résiliation de l'élection de domicile"
$objetPart1 = "Objet"
$objetPart2 = $nbsp + ": "
$objetText = $objetPart1 + $objetPart2 + "Résiliation de l’élection de domicile"
$rng1 = $d.Content
$found1 = $rng1.Find.Execute($objetText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $off1a = $objetPart1.Length
    Split-RunAt $rng1 $off1a

    $off1b = ($objetPart1 + $objetPart2).Length
    Split-RunAt $rng1 $off1b
}

# --- Change 2 -----------------------------------------------------
# Drop the trailing plain space after "...raison suivante<nbsp>: "
$rng2 = $d.Content
$search2 = "suivante" + $nbsp + ": "
$replace2 = "suivante" + $nbsp + ":"
$rng2.Find.Execute($search2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)

# --- Change 3 -----------------------------------------------------
# "...à l'adresse suivante<nbsp>:" (one run) becomes two runs:
# "...à l'adresse " | "suivante<nbsp>:"
$prefix3 = "Nous vous informons que vous pouvez présenter un recours gracieux à l’encontre de cette décision dans les deux mois de sa notification, en adressant votre demande à l’adresse "
$suffix3 = "suivante" + $nbsp + ":"
$fullText3 = $prefix3 + $suffix3
$rng3 = $d.Content
$found3 = $rng3.Find.Execute($fullText3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $off3 = $prefix3.Length
    Split-RunAt $rng3 $off3
}
